# Auto-generated edit script applying scheduled market-price refresh
# to the Gilgamesh Leve-profit tracker workbook (one row per affected Leve).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 91.333336
$ws.Cells.Item(6, 9).Value = 91.333336
$ws.Cells.Item(6, 11).Value = 274.000008
$ws.Cells.Item(6, 13).Value = -162.000008

$ws.Cells.Item(40, 8).Value = 7824.533
$ws.Cells.Item(40, 9).Value = 5416.6665
$ws.Cells.Item(40, 10).Value = 9429.777
$ws.Cells.Item(40, 11).Value = 5416.6665
$ws.Cells.Item(40, 12).Value = 9429.777
$ws.Cells.Item(40, 13).Value = -5241.6665
$ws.Cells.Item(40, 14).Value = -9779.777

$ws.Cells.Item(80, 8).Value = 799.2727
$ws.Cells.Item(80, 9).Value = 243.33333
$ws.Cells.Item(80, 10).Value = 1466.4
$ws.Cells.Item(80, 11).Value = 729.99999
$ws.Cells.Item(80, 12).Value = 4399.200000000001
$ws.Cells.Item(80, 13).Value = 268.00001
$ws.Cells.Item(80, 14).Value = -6395.200000000001

$ws.Cells.Item(83, 8).Value = 799.2727
$ws.Cells.Item(83, 9).Value = 243.33333
$ws.Cells.Item(83, 10).Value = 1466.4
$ws.Cells.Item(83, 11).Value = 2189.99997
$ws.Cells.Item(83, 12).Value = 13197.6
$ws.Cells.Item(83, 13).Value = 2802.00003
$ws.Cells.Item(83, 14).Value = -23181.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1609.3284
$ws.Cells.Item(32, 9).Value = 1609.3284
$ws.Cells.Item(32, 11).Value = 1609.3284
$ws.Cells.Item(32, 13).Value = -1322.3284

$ws.Cells.Item(61, 8).Value = 1634.0968
$ws.Cells.Item(61, 9).Value = 1074.2609
$ws.Cells.Item(61, 10).Value = 3243.625
$ws.Cells.Item(61, 11).Value = 1074.2609
$ws.Cells.Item(61, 12).Value = 3243.625
$ws.Cells.Item(61, 13).Value = -862.2609
$ws.Cells.Item(61, 14).Value = -3667.625

$ws.Cells.Item(74, 8).Value = 2592.077
$ws.Cells.Item(74, 9).Value = 2117.973
$ws.Cells.Item(74, 10).Value = 3761.5334
$ws.Cells.Item(74, 11).Value = 2117.973
$ws.Cells.Item(74, 12).Value = 3761.5334
$ws.Cells.Item(74, 13).Value = -1243.973
$ws.Cells.Item(74, 14).Value = -5509.5334

$ws.Cells.Item(77, 8).Value = 2592.077
$ws.Cells.Item(77, 9).Value = 2117.973
$ws.Cells.Item(77, 10).Value = 3761.5334
$ws.Cells.Item(77, 11).Value = 10589.865
$ws.Cells.Item(77, 12).Value = 18807.667
$ws.Cells.Item(77, 13).Value = -6221.865
$ws.Cells.Item(77, 14).Value = -27543.667

$ws.Cells.Item(136, 8).Value = 1634.0968
$ws.Cells.Item(136, 9).Value = 1074.2609
$ws.Cells.Item(136, 10).Value = 3243.625
$ws.Cells.Item(136, 11).Value = 3222.7827
$ws.Cells.Item(136, 12).Value = 9730.875
$ws.Cells.Item(136, 13).Value = -672.7826999999997
$ws.Cells.Item(136, 14).Value = -14830.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(59, 8).Value = 122998
$ws.Cells.Item(59, 10).Value = 122998
$ws.Cells.Item(59, 12).Value = 122998
$ws.Cells.Item(59, 14).Value = -124692

$ws.Cells.Item(60, 8).Value = 59681.75
$ws.Cells.Item(60, 10).Value = 59681.75
$ws.Cells.Item(60, 12).Value = 59681.75
$ws.Cells.Item(60, 14).Value = -60879.75

$ws.Cells.Item(99, 8).Value = 61145.47
$ws.Cells.Item(99, 9).Value = 68631.53
$ws.Cells.Item(99, 11).Value = 68631.53
$ws.Cells.Item(99, 13).Value = -67133.53

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4012.6
$ws.Cells.Item(31, 9).Value = 3038.9285
$ws.Cells.Item(31, 10).Value = 4661.7144
$ws.Cells.Item(31, 11).Value = 3038.9285
$ws.Cells.Item(31, 12).Value = 4661.7144
$ws.Cells.Item(31, 13).Value = -2743.9285
$ws.Cells.Item(31, 14).Value = -5251.7144

$ws.Cells.Item(34, 8).Value = 4012.6
$ws.Cells.Item(34, 9).Value = 3038.9285
$ws.Cells.Item(34, 10).Value = 4661.7144
$ws.Cells.Item(34, 11).Value = 3038.9285
$ws.Cells.Item(34, 12).Value = 4661.7144
$ws.Cells.Item(34, 13).Value = -2836.9285
$ws.Cells.Item(34, 14).Value = -5065.7144

$ws.Cells.Item(36, 8).Value = 14265.333
$ws.Cells.Item(36, 9).Value = 14500
$ws.Cells.Item(36, 10).Value = 14218.4
$ws.Cells.Item(36, 11).Value = 14500
$ws.Cells.Item(36, 12).Value = 14218.4
$ws.Cells.Item(36, 13).Value = -14112
$ws.Cells.Item(36, 14).Value = -14994.4

$ws.Cells.Item(40, 8).Value = 14265.333
$ws.Cells.Item(40, 9).Value = 14500
$ws.Cells.Item(40, 10).Value = 14218.4
$ws.Cells.Item(40, 11).Value = 14500
$ws.Cells.Item(40, 12).Value = 14218.4
$ws.Cells.Item(40, 13).Value = -14340
$ws.Cells.Item(40, 14).Value = -14538.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(127, 8).Value = 2711.1428
$ws.Cells.Item(127, 10).Value = 2711.1428
$ws.Cells.Item(127, 12).Value = 8133.428400000001
$ws.Cells.Item(127, 14).Value = -18053.4284

$ws.Cells.Item(132, 8).Value = 1904.1666
$ws.Cells.Item(132, 9).Value = 1762
$ws.Cells.Item(132, 10).Value = 1951.5555
$ws.Cells.Item(132, 11).Value = 15858
$ws.Cells.Item(132, 12).Value = 17563.9995
$ws.Cells.Item(132, 13).Value = -13328
$ws.Cells.Item(132, 14).Value = -22623.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 111116440
$ws.Cells.Item(80, 9).Value = 200004800
$ws.Cells.Item(80, 10).Value = 5998.5
$ws.Cells.Item(80, 11).Value = 200004800
$ws.Cells.Item(80, 12).Value = 5998.5
$ws.Cells.Item(80, 13).Value = -200003802
$ws.Cells.Item(80, 14).Value = -7994.5

$ws.Cells.Item(83, 8).Value = 111116440
$ws.Cells.Item(83, 9).Value = 200004800
$ws.Cells.Item(83, 10).Value = 5998.5
$ws.Cells.Item(83, 11).Value = 1000024000
$ws.Cells.Item(83, 12).Value = 29992.5
$ws.Cells.Item(83, 13).Value = -1000019008
$ws.Cells.Item(83, 14).Value = -39976.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 43911.617
$ws.Cells.Item(40, 9).Value = 46737.582
$ws.Cells.Item(40, 11).Value = 46737.582
$ws.Cells.Item(40, 13).Value = -46601.582

$ws.Cells.Item(68, 8).Value = 5143.25
$ws.Cells.Item(68, 9).Value = 5192.6
$ws.Cells.Item(68, 11).Value = 5192.6
$ws.Cells.Item(68, 13).Value = -4443.6

$ws.Cells.Item(71, 8).Value = 5143.25
$ws.Cells.Item(71, 9).Value = 5192.6
$ws.Cells.Item(71, 11).Value = 25963
$ws.Cells.Item(71, 13).Value = -22219

$ws.Cells.Item(92, 8).Value = 49999
$ws.Cells.Item(92, 10).Value = 49999
$ws.Cells.Item(92, 12).Value = 49999
$ws.Cells.Item(92, 14).Value = -54991

$ws.Cells.Item(132, 8).Value = 8642.923000000001
$ws.Cells.Item(132, 9).Value = 5545.3125
$ws.Cells.Item(132, 11).Value = 16635.9375
$ws.Cells.Item(132, 13).Value = -14105.9375

$ws.Cells.Item(136, 8).Value = 3918.4285
$ws.Cells.Item(136, 9).Value = 2001.2916
$ws.Cells.Item(136, 11).Value = 6003.8748
$ws.Cells.Item(136, 13).Value = -3453.8748

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(9, 8).Value = 19968.334
$ws.Cells.Item(9, 9).Value = 19952.5
$ws.Cells.Item(9, 11).Value = 19952.5
$ws.Cells.Item(9, 13).Value = -19812.5

$ws.Cells.Item(14, 8).Value = 7599.75
$ws.Cells.Item(14, 9).Value = 5133
$ws.Cells.Item(14, 11).Value = 5133
$ws.Cells.Item(14, 13).Value = -4965

$ws.Cells.Item(29, 8).Value = 14010
$ws.Cells.Item(29, 9).Value = 14010
$ws.Cells.Item(29, 11).Value = 14010
$ws.Cells.Item(29, 13).Value = -13720

$ws.Cells.Item(126, 8).Value = 2441.125
$ws.Cells.Item(126, 9).Value = 2104.8462
$ws.Cells.Item(126, 11).Value = 6314.5386
$ws.Cells.Item(126, 13).Value = -3844.5386

$ws.Cells.Item(136, 8).Value = 4151.244
$ws.Cells.Item(136, 9).Value = 3691.963
$ws.Cells.Item(136, 11).Value = 11075.889
$ws.Cells.Item(136, 13).Value = -8525.889000000001
